$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.575.67'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '2.494.63'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '2.951.26'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '69.425.08'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.48%  '
$ws.Range('D17').Value = '2.499.42'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '355.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.05%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.49%  '
$ws.Range('D26').Value = '2.621.06'
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').Value = '0.0₃0875'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.63'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '439.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.39%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.71'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.82%  '
$ws.Range('E36').Value = '  -3.30%  '
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('B38').Value = 'POPCAT'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +77.37%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.314'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.80%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '138.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.43%  '
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.507'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.28%  '
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('E51').Value = '  -1.01%  '
